$d = $word.ActiveDocument


$para = $d.Paragraphs(1)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:b/><w:sz w:val="32"/></w:rPr><w:t>18. Recruits</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(14)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:suppressAutoHyphens w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Ashay-Team A</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(15)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:suppressAutoHyphens w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Sidh-Team A</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(18)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:suppressAutoHyphens w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Ral</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(29)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Ashay</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(37)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Ral</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(38)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>Sidh</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$para = $d.Paragraphs(55)
$r = $para.Range
$r.MoveEnd(1, -1)
$r.Text = ""
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>b</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Calibri" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>ar</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')


$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hr = $hdr.Range
$hr.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/header1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.header+xml" pkg:padding="512"><pkg:xmlData><w:hdr xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:p w14:paraId="74EB33FF" w14:textId="53CFB938" w:rsidR="00936FBD" w:rsidRDefault="008967A0"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="7D5C2A4D"><v:line id="_x0000_s2049" style="position:absolute;z-index:-251658752;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9.05pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9.05pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text" from="0,10.8pt" to="466.65pt,10.85pt" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQCGZnnDwQEAAOkDAAAOAAAAZHJzL2Uyb0RvYy54bWysU01v2zAMvQ/YfxB0X5xmaLEYcXpo0V2G&#xA;LdjHD1BkKhagL1BanPz7UYzjdtupw3yQRIp85HuiN/cn78QRMNsYOnmzWEoBQcfehkMnf3x/evdB&#xA;ilxU6JWLATp5hizvt2/fbMbUwioO0fWAgkBCbsfUyaGU1DZN1gN4lRcxQaBLE9GrQiYemh7VSOje&#xA;Navl8q4ZI/YJo4acyft4uZRbxjcGdPliTIYiXCept8Ir8rqva7PdqPaAKg1WT22of+jCKxuo6Az1&#xA;qIoSP9H+BeWtxpijKQsdfRONsRqYA7G5Wf7B5tugEjAXEienWab8/2D15+MOhe3p7aQIytMTVU3G&#xA;lFu6egg7nKycdlgJngz6ulPr4sQ6nmcd4VSEJuftenW7viO59fWueU5MmMtHiF7UQyedDZWiatXx&#xA;Uy5UjEKvIdXtghg7uX5PcNXM0dn+yTrHBh72Dw7FUdXX5a82Twi/hXlboNIgvwu0VXIXOnwqZweX&#xA;Sl/BkBTMiuH1hH+ZFxpoonSdGgajhBpoqJ9X5k4pNRt4TF+ZPydx/RjKnO9tiMgyvGBXj/vYn/k5&#xA;WQCaJ1Zkmv06sC9tlun5D93+AgAA//8DAFBLAwQUAAYACAAAACEA7ourUd4AAAAGAQAADwAAAGRy&#xA;cy9kb3ducmV2LnhtbEyPwU7DMBBE70j8g7VIXCrqpJFCG+JUBQQHDki05e7ESxIRr9PYaQNfz/ZU&#xA;jjszmnmbryfbiSMOvnWkIJ5HIJAqZ1qqFex3L3dLED5oMrpzhAp+0MO6uL7KdWbciT7wuA214BLy&#xA;mVbQhNBnUvqqQav93PVI7H25werA51BLM+gTl9tOLqIolVa3xAuN7vGpwep7O1oFh7T/LH/l+Dhb&#xA;ve3i5X609P78qtTtzbR5ABFwCpcwnPEZHQpmKt1IxotOAT8SFCziFAS7qyRJQJRn4R5kkcv/+MUf&#xA;AAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29u&#xA;dGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAA&#xA;LwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAIZmecPBAQAA6QMAAA4AAAAAAAAAAAAAAAAA&#xA;LgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAO6Lq1HeAAAABgEAAA8AAAAAAAAAAAAA&#xA;AAAAGwQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAAmBQAAAAA=&#xA;" strokeweight=".26mm"><v:stroke joinstyle="miter"/></v:line></w:pict></w:r></w:p></w:hdr></pkg:xmlData></pkg:part></pkg:package>')
